$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 12 ("Enterprises (absolute #)" / 57227) with row 13
# ("Enterprises density (per 1000 people)" / 2) so the density row now
# appears above the absolute-count row. "57227" and "2" are stored as
# text (not numbers), so the swap is done with Copy/PasteSpecial (via an
# unused scratch cell) rather than re-typing the values, which keeps the
# original cell type (text) and style intact instead of Excel
# auto-converting the numeric-looking text into a real number.
$scratchA = $ws.Range("Z100")
$scratchD = $ws.Range("Z101")

$ws.Range("A12").Copy()
$scratchA.PasteSpecial(-4104)   # xlPasteAll
$ws.Range("D12").Copy()
$scratchD.PasteSpecial(-4104)   # xlPasteAll

$ws.Range("A13").Copy()
$ws.Range("A12").PasteSpecial(-4104)   # xlPasteAll
$ws.Range("D13").Copy()
$ws.Range("D12").PasteSpecial(-4104)   # xlPasteAll

$scratchA.Copy()
$ws.Range("A13").PasteSpecial(-4104)   # xlPasteAll
$scratchD.Copy()
$ws.Range("D13").PasteSpecial(-4104)   # xlPasteAll

$scratchA.Clear()
$scratchD.Clear()
$excel.CutCopyMode = $false
